$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "company" (row 4) gets merged into "credits_desc" (row 5) -> delete row 4,
# which shifts credits_desc up to row 4.
$ws.Rows("4").Delete()

# Update the (now-merged) credits_desc value to include both the company
# name and the music credit.
$ws.Range("B4").Value() = "Made by: RENEGADEWARE, Music from: Kevin MacLeod"

# Insert two new rows right after the "summary"/"SUMMARY" row (now row 20)
# for the new "new game" / "continue" UI strings.
$ws.Rows("21:22").Insert()
$ws.Range("A21").Value() = "new"
$ws.Range("B21").Value() = "NEW GAME"
$ws.Range("A22").Value() = "continue"
$ws.Range("B22").Value() = "CONTINUE"

# Restore the selection to B4, matching the saved workbook state.
$ws.Range("B4").Select()
